$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Office Visit" answer (Contact category): UAE -> DUBAI
$ws.Range("D28").Value = "We are located at Entrance no:4 First Floor Office no 2056-A, Dubai Cargo Village, DUBAI."

# Update the "Submit Feedback" answer (Feedback category): feedback@ukfservices.com -> info@ukfservices.com
$ws.Range("D29").Value = "You can provide feedback via our website form or email info@ukfservices.com."

# Reflect the updated viewport/zoom/selection left behind by the author's Excel session
$excel.ActiveWindow.Zoom = 115
$ws.Range("D34").Select() | Out-Null
